$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Observación" -> "Observaciones" (Interna #5 row, "Observación." bullet)
# ---------------------------------------------------------------------------
$obsRng = $d.Content
$obsRng.Find.Execute("Observación", $true, $false, $false, $false, $false, $true, 1, $false, "Observaciones", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Give the "Entrega de producto" paragraph (Externa #4 row) the same
#    bullet numbering (numId 35) already used by the other bulleted rows,
#    and move the "_GoBack" bookmark from the trailing empty paragraph onto
#    the start of this paragraph.
# ---------------------------------------------------------------------------

# Locate an existing paragraph that already carries numId 35 so we can reuse
# its exact ListTemplate (ContinuePreviousList keeps the same numId instead
# of minting a new list definition).
$srcRng = $d.Content
$srcRng.Find.Execute("Revisión de avances de producto", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$allParas = $d.Paragraphs
$srcPara = $null
for ($i = 1; $i -le $allParas.Count; $i++) {
    $p = $allParas.Item($i)
    if ($p.Range.Start -le $srcRng.Start -and $p.Range.End -ge $srcRng.End) {
        $srcPara = $p
        break
    }
}
$listTemplate = $srcPara.Range.ListFormat.ListTemplate

# Locate the target paragraph ("Entrega de producto").
$tgtRng = $d.Content
$tgtRng.Find.Execute("Entrega de producto", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tgtPara = $null
for ($i = 1; $i -le $allParas.Count; $i++) {
    $p = $allParas.Item($i)
    if ($p.Range.Start -le $tgtRng.Start -and $p.Range.End -ge $tgtRng.End) {
        $tgtPara = $p
        break
    }
}

# Apply the list numbering, continuing the existing list (numId 35).
$tgtPara.Range.ListFormat.ApplyListTemplate($listTemplate, $true)

# Move the "_GoBack" bookmark onto the start of the target paragraph.
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$tgtStartRng = $d.Content
$tgtStartRng.Find.Execute("Entrega de producto", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$newBookmarkRange = $d.Range($tgtStartRng.Start, $tgtStartRng.Start)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange) | Out-Null
